# Commit: "Updated notebook, reran simulation"
#
# The underlying averaging notebook picked up two new measurement techniques
# ("Holden" and "Rizzie Spiral", inserted right after "Spiral5") and renamed
# "Thomas Hex" to "Matthies Hex". Re-running the simulation also appended two more
# rows of results ("Michael-CCHex", "Michael-SNHex") at the bottom of the sheet.
# Column A keeps the running 0-based index used elsewhere in the notebook, and the
# data columns (C:W) keep reporting 1 for every HKL/pair combination, unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row labels for B4:B29 (technique names shift down the list) ---
$labels = New-Object 'object[,]' 26,1
$labels[0,0] = "Holden"
$labels[1,0] = "Rizzie Spiral"
$labels[2,0] = "RotRing OmegaMax-90"
$labels[3,0] = "Equal Angle"
$labels[4,0] = "Tilt Rotate"
$labels[5,0] = "CLR"
$labels[6,0] = "Rizzie Hex"
$labels[7,0] = "Matthies Hex"
$labels[8,0] = "Tilt Rotate_Partial"
$labels[9,0] = "RotRing OmegaMax-60"
$labels[10,0] = "Equal Angle_Partial"
$labels[11,0] = "Rizzie Hex_Partial"
$labels[12,0] = "ND Single"
$labels[13,0] = "RD Single"
$labels[14,0] = "TD Single"
$labels[15,0] = "Morris Single"
$labels[16,0] = "Ring Perpendicular to ND"
$labels[17,0] = "Ring Perpendicular to RD"
$labels[18,0] = "Ring Perpendicular to TD"
$labels[19,0] = "OffsetFTD"
$labels[20,0] = "OffsetATD"
$labels[21,0] = "OffsetF45"
$labels[22,0] = "OffsetA45"
$labels[23,0] = "OffsetFRD"
$labels[24,0] = "OffsetARD"
$labels[25,0] = "Gaussian Quadrature"
$ws.Range("B4:B29").Value = $labels

# --- Append two new simulation result rows (30 and 31) ---
$newRows = New-Object 'object[,]' 2,23
$newRows[0,0] = 28
$newRows[0,1] = "Michael-CCHex"
$newRows[1,0] = 29
$newRows[1,1] = "Michael-SNHex"
for ($col = 2; $col -lt 23; $col++) {
    $newRows[0,$col] = 1
    $newRows[1,$col] = 1
}
$ws.Range("A30:W31").Value = $newRows

# Match column A formatting (bold, centered, bordered index style) used by the rest of the table
$ws.Range("A29").Copy() | Out-Null
$ws.Range("A30:A31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A1").Select() | Out-Null

Write-Host "Used range after edit: $($ws.UsedRange.Address())"
